$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.788428902626038
$ws.Range("B1").Value = 1.892302393913269
$ws.Range("C1").Value = 2.293514728546143
$ws.Range("D1").Value = 2.158591270446777
$ws.Range("E1").Value = 3.037123680114746
